# Applies the 2025-12-04 18:29 JST re-scrape: refreshes the timestamp on the
# existing rows, inserts two new postings above the previous #1 item, shifts
# the rest down, appends one new posting at the bottom, widens column B, and
# re-links the URL hyperlinks so they keep pointing at the right row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rows 2-13: only the fetched-at timestamp (column A) changes.
$ws.Range("A2").Value = '2025-12-04 18:29:49'
$ws.Range("A3").Value = '2025-12-04 18:29:49'
$ws.Range("A4").Value = '2025-12-04 18:29:49'
$ws.Range("A5").Value = '2025-12-04 18:29:49'
$ws.Range("A6").Value = '2025-12-04 18:29:49'
$ws.Range("A7").Value = '2025-12-04 18:29:49'
$ws.Range("A8").Value = '2025-12-04 18:29:49'
$ws.Range("A9").Value = '2025-12-04 18:29:49'
$ws.Range("A10").Value = '2025-12-04 18:29:49'
$ws.Range("A11").Value = '2025-12-04 18:29:49'
$ws.Range("A12").Value = '2025-12-04 18:29:49'
$ws.Range("A13").Value = '2025-12-04 18:29:49'

# 2) Drop all existing URL hyperlinks up front; rows 14-18 are about to be
#    overwritten/relocated, so stale hyperlink-to-row mappings would
#    otherwise survive. They are re-created in step 4 below.
$ws.Hyperlinks.Delete()

# 3) Rewrite rows 14-21 in full with the refreshed listings (two new rows
#    inserted above the old #1 item, the old rows 14-18 shifted down to
#    16-20, and one brand-new row 21 appended).
# row 14
$ws.Range("A14").Value = '2025-12-04 18:29:49'
$ws.Range("B14").Value = '初回 【急募】ECサイトの要件定義や基本設計ができる方を募集(1人月、フルリモート可、2025年12月〜)'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5425629'
$ws.Range("G14").Value = 45
$ws.Range("H14").Value = ""

# row 15
$ws.Range("A15").Value = '2025-12-04 18:29:49'
$ws.Range("B15").Value = '【UTAGE構築代行】各種初期設定・ステップ配信・会員サイトの作成など'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5447344'
$ws.Range("G15").Value = 38
$ws.Range("H15").Value = '◇サイト'

# row 16
$ws.Range("A16").Value = '2025-12-04 18:29:49'
$ws.Range("B16").Value = '【急募】宝くじ予想サイトの構築'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5446997'
$ws.Range("G16").Value = 38
$ws.Range("H16").Value = '◇サイト'

# row 17
$ws.Range("A17").Value = '2025-12-04 18:29:49'
$ws.Range("B17").Value = '【急募】住所リストから太陽光パネル設置を自動判定するシステム'
$ws.Range("C17").Value = 'システム開発'
$ws.Range("D17").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E17").Value = '期限情報なし'
$ws.Range("F17").Value = 'https://www.lancers.jp/work/detail/5447102'
$ws.Range("G17").Value = 33

# row 18
$ws.Range("A18").Value = '2025-12-04 18:29:49'
$ws.Range("B18").Value = '【SESエンジニア募集】多様なプロジェクトに参画可能!'
$ws.Range("C18").Value = 'システム開発'
$ws.Range("D18").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E18").Value = '期限情報なし'
$ws.Range("F18").Value = 'https://www.lancers.jp/work/detail/5437544'
$ws.Range("G18").Value = 25

# row 19
$ws.Range("A19").Value = '2025-12-04 18:29:49'
$ws.Range("B19").Value = '【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え'
$ws.Range("C19").Value = 'システム開発'
$ws.Range("D19").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E19").Value = '期限情報なし'
$ws.Range("F19").Value = 'https://www.lancers.jp/work/detail/5443568'
$ws.Range("G19").Value = 13

# row 20
$ws.Range("A20").Value = '2025-12-04 18:29:49'
$ws.Range("B20").Value = '注目 【電子工作】蓋を開くとmp3再生するスピーカー制作依頼'
$ws.Range("C20").Value = 'システム開発'
$ws.Range("D20").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E20").Value = '期限情報なし'
$ws.Range("F20").Value = 'https://www.lancers.jp/work/detail/5446806'
$ws.Range("G20").Value = 13

# row 21
$ws.Range("A21").Value = '2025-12-04 18:29:49'
$ws.Range("B21").Value = '限定公開 限定公開の仕事'
$ws.Range("C21").Value = 'システム開発'
$ws.Range("D21").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E21").Value = '期限情報なし'
$ws.Range("F21").Value = 'https://www.lancers.jp/work/detail/5447357'
$ws.Range("G21").Value = 10

# 4) Re-add one hyperlink per data row (F2:F21) targeting the URL text now
#    in each cell.
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5423720') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5446833') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5419380') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5447137') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5446990') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5446867') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5441557') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5447021') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5441568') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5446668') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5431107') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5446849') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5425629') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5447344') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5446997') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5447102') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F18"), 'https://www.lancers.jp/work/detail/5437544') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F19"), 'https://www.lancers.jp/work/detail/5443568') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F20"), 'https://www.lancers.jp/work/detail/5446806') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F21"), 'https://www.lancers.jp/work/detail/5447357') | Out-Null

# 5) Column B (title) widens from 49 to 55 characters. ColumnWidth adds
#    Excel's fixed ~0.8333-character padding on top of the assigned value,
#    so back it out to land exactly on 55 in the saved file.
$ws.Columns.Item(2).ColumnWidth = 54.16666666666667
